$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 156 (shifts the existing rows 156..237 down to 157..238)
$ws.Rows.Item(156).Insert()

$ws.Range("A156").Value = 4
$ws.Range("B156").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C156").Value = "Los Lagos"
$ws.Range("D156").Value = 44572
$ws.Range("E156").Value = 10
$ws.Range("F156").Value = 100112045
$ws.Range("G156").Value = "Zapallo"
$ws.Range("H156").Value = "Paine"
$ws.Range("I156").Value = "1a nueva(o)"
$ws.Range("J156").Value = 1100
$ws.Range("K156").Value = 500
$ws.Range("L156").Value = 550
$ws.Range("M156").Value = 527
$ws.Range("N156").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O156").Value = "Región de O'Higgins"
$ws.Range("P156").Value = 527
$ws.Range("Q156").Value = 1
$ws.Range("R156").Value = "Hortaliza"
